$d = $word.ActiveDocument

# 1. Update the date heading paragraph
$d.Content.Find.Execute("2025-07-13 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-07-14 Monday", 2) | Out-Null

# 2. Update the division-fact table cells.
# Assigning Cell.Range.Text replaces only the cell's text while keeping the
# existing run/paragraph formatting (font, size, alignment) intact.
$t = $d.Tables.Item(1)

$rows = @(1, 5, 9, 13, 17)
$newValues = @(
    @("63÷7=9, 0",  "63÷9=7, 0",  "31÷7=4, 3",  "20÷5=4, 0",  "82÷6=13, 4"),
    @("11÷9=1, 2",  "24÷2=12, 0", "36÷8=4, 4",  "78÷5=15, 3", "64÷3=21, 1"),
    @("78÷8=9, 6",  "19÷6=3, 1",  "65÷9=7, 2",  "20÷7=2, 6",  "17÷7=2, 3"),
    @("17÷3=5, 2",  "64÷2=32, 0", "33÷5=6, 3",  "61÷2=30, 1", "50÷3=16, 2"),
    @("39÷5=7, 4",  "32÷4=8, 0",  "78÷6=13, 0", "47÷3=15, 2", "77÷5=15, 2")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowIndex = $rows[$i]
    $values = $newValues[$i]
    for ($col = 1; $col -le $values.Length; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
